# Update "想去人数" (attendance count) values on three sheets:
#   展览 (Exhibition)      -> sheet1
#   本地生活 (Local Life)  -> sheet3
#   全部类型 (All Types)   -> sheet4
# "演出" (Performance) is unchanged.

$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value  = 874
$ws.Range("F7").Value  = 458
$ws.Range("F9").Value  = 2143
$ws.Range("F11").Value = 278
$ws.Range("F12").Value = 115
$ws.Range("F13").Value = 1046
$ws.Range("F14").Value = 174
$ws.Range("F15").Value = 2176
$ws.Range("F16").Value = 640
$ws.Range("F17").Value = 11852
$ws.Range("F18").Value = 1209
$ws.Range("F21").Value = 10
$ws.Range("F24").Value = 259
$ws.Range("F27").Value = 13

# --- 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 474

# --- 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 474
$ws.Range("F10").Value = 874
$ws.Range("F12").Value = 458
$ws.Range("F14").Value = 2143
$ws.Range("F16").Value = 278
$ws.Range("F18").Value = 115
$ws.Range("F20").Value = 1047
$ws.Range("F22").Value = 174
$ws.Range("F25").Value = 2176
$ws.Range("F26").Value = 640
$ws.Range("F29").Value = 1209
$ws.Range("F32").Value = 10
$ws.Range("F38").Value = 259
$ws.Range("F49").Value = 13

$wb.Save()
